# Rename "Acme mAb N" -> "VD-Crotty N" on the Antibodies sheet (rows 2-11)
$wb = $excel.ActiveWorkbook
$wsAntibodies = $wb.Worksheets.Item("Antibodies")
$wsTerminology = $wb.Worksheets.Item("Terminology")

for ($i = 1; $i -le 10; $i++) {
    $row = $i + 1
    $wsAntibodies.Range("A$row").Value = "VD-Crotty $i"
}

# Remove rows 16 and 17 (kappa, lambda) from the Terminology sheet
$wsTerminology.Rows.Item(16).Delete()
$wsTerminology.Rows.Item(16).Delete()

# Update the data validation formula for column C on the Antibodies sheet
# to reference the shrunk Terminology!B range (B2:B15 instead of B2:B17)
$rangeC = $wsAntibodies.Range("C2:C100")
$rangeC.Validation.Formula1 = "=Terminology!`$B`$2:`$B`$15"
